$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 0
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0
$ws.Range("C8").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 0
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 0
$ws.Range("C11").Value = 1
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 0
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 0
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 0
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 0
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 0
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 0
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 0
$ws.Range("C21").Value = 1
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 0
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 0
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 0
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 0
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 0
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 0
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 0
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 0
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 0
$ws.Range("C32").Value = 1
$ws.Range("D32").Value = 0
$ws.Range("C33").Value = 1
$ws.Range("C34").Value = 1
$ws.Range("C35").Value = 1
$ws.Range("D35").Value = 0
$ws.Range("C36").Value = 1
$ws.Range("D36").Value = 0
$ws.Range("C37").Value = 1
$ws.Range("D37").Value = 0
$ws.Range("C38").Value = 1
$ws.Range("C39").Value = 1
$ws.Range("D39").Value = 0
$ws.Range("C40").Value = 1
$ws.Range("D40").Value = 0
$ws.Range("C41").Value = 1
$ws.Range("C42").Value = 1
$ws.Range("D42").Value = 0
$ws.Range("C43").Value = 1
$ws.Range("D43").Value = 0
$ws.Range("C44").Value = 1
$ws.Range("D44").Value = 0
$ws.Range("C45").Value = 1
$ws.Range("D45").Value = 0
$ws.Range("C46").Value = 1
$ws.Range("D46").Value = 0
$ws.Range("C47").Value = 1
$ws.Range("D47").Value = 0
$ws.Range("C48").Value = 1
$ws.Range("D48").Value = 0
$ws.Range("C49").Value = 1
$ws.Range("D49").Value = 0
$ws.Range("C50").Value = 1
$ws.Range("D50").Value = 0
$ws.Range("C51").Value = 1
$ws.Range("C52").Value = 1
$ws.Range("D52").Value = 0
$ws.Range("C53").Value = 1
$ws.Range("D53").Value = 0
$ws.Range("C54").Value = 1
$ws.Range("D54").Value = 0
$ws.Range("C55").Value = 1
$ws.Range("D55").Value = 0
$ws.Range("C56").Value = 1
$ws.Range("D56").Value = 0
$ws.Range("C57").Value = 1
$ws.Range("D57").Value = 0
$ws.Range("C58").Value = 1
$ws.Range("D58").Value = 0
$ws.Range("C59").Value = 1
$ws.Range("C61").Value = 1
$ws.Range("C62").Value = 1
$ws.Range("C63").Value = 1
$ws.Range("D63").Value = 0
$ws.Range("C64").Value = 1
$ws.Range("D64").Value = 0
$ws.Range("C65").Value = 1
$ws.Range("C66").Value = 1
$ws.Range("D66").Value = 0
$ws.Range("C67").Value = 1
$ws.Range("D67").Value = 0
$ws.Range("C68").Value = 1
$ws.Range("D68").Value = 0
$ws.Range("C69").Value = 1
$ws.Range("D69").Value = 0
$ws.Range("C71").Value = 1
$ws.Range("C72").Value = 1
$ws.Range("D72").Value = 0
$ws.Range("C73").Value = 1
$ws.Range("D73").Value = 0
$ws.Range("C74").Value = 1
$ws.Range("D74").Value = 0
$ws.Range("C75").Value = 1
$ws.Range("C76").Value = 1
$ws.Range("D76").Value = 0
$ws.Range("C77").Value = 1
$ws.Range("C78").Value = 1
$ws.Range("D78").Value = 0
$ws.Range("C79").Value = 1
$ws.Range("D79").Value = 0
$ws.Range("C80").Value = 1
$ws.Range("D80").Value = 0
$ws.Range("C81").Value = 1
$ws.Range("D81").Value = 0
$ws.Range("C82").Value = 1
$ws.Range("D82").Value = 0
$ws.Range("C83").Value = 1
$ws.Range("D83").Value = 0
$ws.Range("C84").Value = 1
$ws.Range("D84").Value = 0
$ws.Range("C85").Value = 1
$ws.Range("D85").Value = 0
$ws.Range("C86").Value = 1
$ws.Range("D86").Value = 0
$ws.Range("C88").Value = 1
$ws.Range("C90").Value = 1
$ws.Range("D90").Value = 0
$ws.Range("D91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 0
$ws.Range("C96").Value = 1
$ws.Range("D96").Value = 0
$ws.Range("C97").Value = 1
$ws.Range("C98").Value = 1
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("C101").Value = 1
$ws.Range("D101").Value = 0
$ws.Range("C102").Value = 1
$ws.Range("D102").Value = 0
$ws.Range("C103").Value = 1
$ws.Range("D103").Value = 0
$ws.Range("D104").Value = 0
$ws.Range("C106").Value = 1
$ws.Range("D106").Value = 0
$ws.Range("D107").Value = 0
$ws.Range("C108").Value = 1
$ws.Range("D108").Value = 0
$ws.Range("D110").Value = 0
$ws.Range("D111").Value = 0
$ws.Range("C112").Value = 1
$ws.Range("D112").Value = 0
$ws.Range("D114").Value = 0
$ws.Range("D115").Value = 0
$ws.Range("C117").Value = 1
$ws.Range("D117").Value = 0
$ws.Range("D118").Value = 0
$ws.Range("D119").Value = 0
$ws.Range("C121").Value = 1
$ws.Range("D121").Value = 0
$ws.Range("D122").Value = 0
$ws.Range("D123").Value = 0
$ws.Range("D124").Value = 0
$ws.Range("D125").Value = 0
$ws.Range("D127").Value = 0
$ws.Range("D129").Value = 0
$ws.Range("C130").Value = 1
$ws.Range("D130").Value = 0
$ws.Range("D131").Value = 0
$ws.Range("C132").Value = 1
$ws.Range("D132").Value = 0
$ws.Range("D133").Value = 0
$ws.Range("C134").Value = 1
$ws.Range("D134").Value = 0
$ws.Range("D135").Value = 0
$ws.Range("D136").Value = 0
$ws.Range("D137").Value = 0
$ws.Range("D139").Value = 0
$ws.Range("D141").Value = 0
$ws.Range("C143").Value = 1
$ws.Range("D143").Value = 0
$ws.Range("C144").Value = 1
$ws.Range("D144").Value = 0
$ws.Range("C145").Value = 1
$ws.Range("D145").Value = 0
$ws.Range("D146").Value = 0
$ws.Range("D149").Value = 0
$ws.Range("D150").Value = 0
